$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 20,14
$data[0,0] = 0.026992
$data[0,1] = 0.08097599999999999
$data[0,2] = 0.004182906599909731
$data[0,3] = 0.00420788870005516
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 1.525469666666667
$data[0,7] = 4.576409
$data[0,8] = 0.02884292053092702
$data[0,9] = 0.03941949540871108
$data[0,10] = 0.04117547724266667
$data[0,11] = 0.370579295184
$data[0,12] = 0.0001206472426494865
$data[0,13] = 0.0001658728492921916
$data[1,0] = 0.026992
$data[1,1] = 0.08097599999999999
$data[1,2] = 0.004182906599909731
$data[1,3] = 0.00420788870005516
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 6.821209666666667
$data[1,7] = 20.463629
$data[1,8] = 0.1289724814852373
$data[1,9] = 0.1762661356122381
$data[1,10] = 0.1841180913226667
$data[1,11] = 1.657062821904
$data[1,12] = 0.0005394798440113347
$data[1,13] = 0.0007417082802451271
$data[2,0] = 0.026992
$data[2,1] = 0.08097599999999999
$data[2,2] = 0.004182906599909731
$data[2,3] = 0.00420788870005516
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 1.197510666666667
$data[2,7] = 3.592532
$data[2,8] = 0.02264201363575945
$data[2,9] = 0.03094474263109954
$data[2,10] = 0.03232320791466667
$data[2,11] = 0.290908871232
$data[2,12] = 0.00009470942827226434
$data[2,13] = 0.000130212032843519
$data[3,0] = 0.026992
$data[3,1] = 0.08097599999999999
$data[3,2] = 0.004182906599909731
$data[3,3] = 0.00420788870005516
$data[3,4] = 2
$data[3,5] = 1
$data[3,6] = 42.571562
$data[3,7] = 85.143124
$data[3,8] = 0.804924677608644
$data[3,9] = 0.7333913960927264
$data[3,10] = 1.149091601504
$data[3,11] = 6.894549609024
$data[3,12] = 0.003366924746399409
$data[3,13] = 0.003086029368336262
$data[4,0] = 0.026992
$data[4,1] = 0.08097599999999999
$data[4,2] = 0.004182906599909731
$data[4,3] = 0.00420788870005516
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.7731246666666668
$data[4,7] = 2.319374
$data[4,8] = 0.01461790673943223
$data[4,9] = 0.01997823025522497
$data[4,10] = 0.02086818100266667
$data[4,11] = 0.187813629024
$data[4,12] = 0.00006114533857723602
$data[4,13] = 0.00008406616933806128
$data[5,0] = 6.300519666666666
$data[5,1] = 18.901559
$data[5,2] = 0.976381346197431
$data[5,3] = 0.9822127115383066
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.525469666666667
$data[5,7] = 4.576409
$data[5,8] = 0.02884292053092702
$data[5,9] = 0.03941949540871108
$data[5,10] = 9.611251635736778
$data[5,11] = 86.50126472163099
$data[5,12] = 0.02816168957625204
$data[5,13] = 0.03871832947286193
$data[6,0] = 6.300519666666666
$data[6,1] = 18.901559
$data[6,2] = 0.976381346197431
$data[6,3] = 0.9822127115383066
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 6.821209666666667
$data[6,7] = 20.463629
$data[6,8] = 0.1289724814852373
$data[6,9] = 0.1762661356122381
$data[6,10] = 42.97716565529011
$data[6,11] = 386.794490897611
$data[6,12] = 0.1259263250949793
$data[6,13] = 0.1731308390120752
$data[7,0] = 6.300519666666666
$data[7,1] = 18.901559
$data[7,2] = 0.976381346197431
$data[7,3] = 0.9822127115383066
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.197510666666667
$data[7,7] = 3.592532
$data[7,8] = 0.02264201363575945
$data[7,9] = 0.03094474263109954
$data[7,10] = 7.544939506376445
$data[7,11] = 67.904455557388
$data[7,12] = 0.02210723975430341
$data[7,13] = 0.03039431956754731
$data[8,0] = 6.300519666666666
$data[8,1] = 18.901559
$data[8,2] = 0.976381346197431
$data[8,3] = 0.9822127115383066
$data[8,4] = 2
$data[8,5] = 1
$data[8,6] = 42.571562
$data[8,7] = 85.143124
$data[8,8] = 0.804924677608644
$data[8,9] = 0.7333913960927264
$data[8,10] = 268.2229636217193
$data[8,11] = 1609.337781730316
$data[8,12] = 0.7859134403110609
$data[8,13] = 0.720346351775101
$data[9,0] = 6.300519666666666
$data[9,1] = 18.901559
$data[9,2] = 0.976381346197431
$data[9,3] = 0.9822127115383066
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.7731246666666668
$data[9,7] = 2.319374
$data[9,8] = 0.01461790673943223
$data[9,9] = 0.01997823025522497
$data[9,10] = 4.871087167118445
$data[9,11] = 43.839784504066
$data[9,12] = 0.01427265146083534
$data[9,13] = 0.01962287171072115
$data[10,0] = 0.1149325
$data[10,1] = 0.229865
$data[10,2] = 0.01781090370458377
$data[10,3] = 0.01194485200600399
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 1.525469666666667
$data[10,7] = 4.576409
$data[10,8] = 0.02884292053092702
$data[10,9] = 0.03941949540871108
$data[10,10] = 0.1753260424641667
$data[10,11] = 1.051956254785
$data[10,12] = 0.0005137184801353032
$data[10,13] = 0.0004708600388084077
$data[11,0] = 0.1149325
$data[11,1] = 0.229865
$data[11,2] = 0.01781090370458377
$data[11,3] = 0.01194485200600399
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 6.821209666666667
$data[11,7] = 20.463629
$data[11,8] = 0.1289724814852373
$data[11,9] = 0.1762661356122381
$data[11,10] = 0.7839786800141667
$data[11,11] = 4.703872080085
$data[11,12] = 0.002297116448274775
$data[11,13] = 0.002105472903558414
$data[12,0] = 0.1149325
$data[12,1] = 0.229865
$data[12,2] = 0.01781090370458377
$data[12,3] = 0.01194485200600399
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 1.197510666666667
$data[12,7] = 3.592532
$data[12,8] = 0.02264201363575945
$data[12,9] = 0.03094474263109954
$data[12,10] = 0.1376328946966667
$data[12,11] = 0.82579736818
$data[12,12] = 0.0004032747245443843
$data[12,13] = 0.0003696303710923667
$data[13,0] = 0.1149325
$data[13,1] = 0.229865
$data[13,2] = 0.01781090370458377
$data[13,3] = 0.01194485200600399
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 42.571562
$data[13,7] = 85.143124
$data[13,8] = 0.804924677608644
$data[13,9] = 0.7333913960927264
$data[13,10] = 4.892856049564999
$data[13,11] = 19.57142419826
$data[13,12] = 0.01433643592233069
$data[13,13] = 0.008760251688804272
$data[14,0] = 0.1149325
$data[14,1] = 0.229865
$data[14,2] = 0.01781090370458377
$data[14,3] = 0.01194485200600399
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 0.7731246666666668
$data[14,7] = 2.319374
$data[14,8] = 0.01461790673943223
$data[14,9] = 0.01997823025522497
$data[14,10] = 0.08885715075166667
$data[14,11] = 0.53314290451
$data[14,12] = 0.0002603581292986136
$data[14,13] = 0.0002386370037405336
$data[15,0] = 0.010485
$data[15,1] = 0.031455
$data[15,2] = 0.001624843498075486
$data[15,3] = 0.001634547755634201
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 1.525469666666667
$data[15,7] = 4.576409
$data[15,8] = 0.02884292053092702
$data[15,9] = 0.03941949540871108
$data[15,10] = 0.015994549455
$data[15,11] = 0.143950945095
$data[15,12] = 0.00004686523189018472
$data[15,13] = 0.00006443304774854138
$data[16,0] = 0.010485
$data[16,1] = 0.031455
$data[16,2] = 0.001624843498075486
$data[16,3] = 0.001634547755634201
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 6.821209666666667
$data[16,7] = 20.463629
$data[16,8] = 0.1289724814852373
$data[16,9] = 0.1762661356122381
$data[16,10] = 0.07152038335500001
$data[16,11] = 0.643683450195
$data[16,12] = 0.0002095600979719489
$data[16,13] = 0.0002881154163592974
$data[17,0] = 0.010485
$data[17,1] = 0.031455
$data[17,2] = 0.001624843498075486
$data[17,3] = 0.001634547755634201
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 1.197510666666667
$data[17,7] = 3.592532
$data[17,8] = 0.02264201363575945
$data[17,9] = 0.03094474263109954
$data[17,10] = 0.01255589934
$data[17,11] = 0.11300309406
$data[17,12] = 0.00003678972863940026
$data[17,13] = 0.00005058065961634173
$data[18,0] = 0.010485
$data[18,1] = 0.031455
$data[18,2] = 0.001624843498075486
$data[18,3] = 0.001634547755634201
$data[18,4] = 2
$data[18,5] = 1
$data[18,6] = 42.571562
$data[18,7] = 85.143124
$data[18,8] = 0.804924677608644
$data[18,9] = 0.7333913960927264
$data[18,10] = 0.44636282757
$data[18,11] = 2.67817696542
$data[18,12] = 0.001307876628852912
$data[18,13] = 0.001198763260484799
$data[19,0] = 0.010485
$data[19,1] = 0.031455
$data[19,2] = 0.001624843498075486
$data[19,3] = 0.001634547755634201
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 0.7731246666666668
$data[19,7] = 2.319374
$data[19,8] = 0.01461790673943223
$data[19,9] = 0.01997823025522497
$data[19,10] = 0.00810621213
$data[19,11] = 0.07295590917
$data[19,12] = 0.0000237518107210403
$data[19,13] = 0.00003265537142522126

$ws.Range("G2:T21").Value = $data

Write-Host "Updated G2:T21"
